# "Finished Product Testing Sequence"
# - Mark every product row (2-40) on the "Products" sheet as tested
#   (new boolean "Test Result" column E = TRUE).
# - The "Products" sheet becomes the active/selected tab (was "Orders"),
#   scrolled so row 31 is at the top.

$wb = $excel.ActiveWorkbook

$wsProducts = $wb.Worksheets.Item("Products")

# Fill in the "Test Result" column (E) with TRUE for every data row.
$wsProducts.Range("E2:E40").Value = $true

# Make "Products" the active sheet/tab (this also clears tabSelected on
# whichever sheet was previously active, i.e. "Orders").
$wsProducts.Activate()

# Scroll the view so row 31 is at the top-left of the window, and restore
# the original selected cell for this sheet.
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$wsProducts.Range("E1").Select() | Out-Null
